# Apply "Updated symbol list" edits to the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Simple price / label updates (single-cell value replacements)
Set-TextValue "D3"  "22.22"
Set-TextValue "D4"  "5.364"
Set-TextValue "D5"  "0.05873"
Set-TextValue "D6"  "3.387"
Set-TextValue "D7"  "6.381"
Set-TextValue "D8"  "0.8123"
Set-TextValue "D9"  "0.9568"
Set-TextValue "D10" "0.1422"
Set-TextValue "D11" "0.03531"
Set-TextValue "D12" "0.07375"
Set-TextValue "D13" "0.03034"
Set-TextValue "D14" "4.418"
Set-TextValue "D15" "0.09397"
Set-TextValue "D16" "0.001590"
Set-TextValue "D17" "0.04813"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006119"
Set-TextValue "D20" "0.004090"
Set-TextValue "D21" "0.0009842"
Set-TextValue "D22" "0.00009701"
Set-TextValue "D23" "3.687"
Set-TextValue "D27" "0.0002472"
Set-TextValue "D40" "0.03859"

# Rows 41-43 get reshuffled (coins rotate: BKEXToken/CEJI/KickToken -> KickToken/BKEXToken/CEJI)
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006641"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1075"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003000"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D44" "0.005761"
Set-TextValue "D45" "0.00005670"

Set-TextValue "D48" "0.07578"
Set-TextValue "E48" "47BOLOBOLO"
